$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Expand the two "comparison label" text values (the real content change)
# Order matters for shared-string table layout: update A5 (computer) first,
# then A3 (eetet), so new/changed strings are appended in the same order
# Excel produced them.
$ws.Range("A5").Value = "computer computer computer"
$ws.Range("A3").Value = "eetet eetet eetet eetet eetet"

# Update the active selection to A5, as recorded in the saved view state
$ws.Range("A5").Select()
